$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds one daily price record per row (rows 2-186). This edit
# inserts a brand-new record as row 60, which pushes every record that used
# to live in rows 60-186 down by one row (to rows 61-187). The row-invariant
# identity columns (A Mercado ID, B Mercado, C Region, E Codreg, F Categoria ID,
# G Categoria, H Variedad, R Clasificacion) are identical on every data row in
# this sheet, so only the "payload" columns D,I,J,K,L,M,N,O,P,Q actually need
# to move.

# 1) Snapshot the payload columns (D:Q) for the rows that will shift,
#    BEFORE writing anything, so later writes can't clobber data we still
#    need to read.
$old = $ws.Range("D60:Q186").Value2

# 2) Write that snapshot one row lower (D61:Q187) - this is the "shift down
#    by one" that turns old row 60 into new row 61, ..., old row 186 into
#    new row 187.
$ws.Range("D61:Q187").Value2 = $old

# 3) The new row 187 is a genuinely new row, so its identity columns
#    (A,B,C,E,F,G,H,R) need to be populated - copy them from the row above
#    (186), which (like every other row) carries the same constant values.
$ws.Cells.Item(187, 1).Value2  = $ws.Cells.Item(186, 1).Value2   # A Mercado ID
$ws.Cells.Item(187, 2).Value2  = $ws.Cells.Item(186, 2).Value2   # B Mercado
$ws.Cells.Item(187, 3).Value2  = $ws.Cells.Item(186, 3).Value2   # C Region
$ws.Cells.Item(187, 5).Value2  = $ws.Cells.Item(186, 5).Value2   # E Codreg
$ws.Cells.Item(187, 6).Value2  = $ws.Cells.Item(186, 6).Value2   # F Categoria ID
$ws.Cells.Item(187, 7).Value2  = $ws.Cells.Item(186, 7).Value2   # G Categoria
$ws.Cells.Item(187, 8).Value2  = $ws.Cells.Item(186, 8).Value2   # H Variedad
$ws.Cells.Item(187, 18).Value2 = $ws.Cells.Item(186, 18).Value2  # R Clasificacion

# The Fecha column (D) carries a date/time number format - make sure the
# newly created row 187 cell gets the same display format as the rest of
# column D instead of defaulting to General.
$ws.Cells.Item(187, 4).NumberFormat = $ws.Cells.Item(186, 4).NumberFormat

# 4) Finally, overwrite row 60 with the brand-new record's values. Its
#    identity columns (A,B,C,E,F,G,H,R) and the untouched payload columns
#    (I Calidad, N Unidad de comercializacion, O Origen, Q Kg o Unidades)
#    already hold the right values and are left as-is.
$ws.Cells.Item(60, 4).Value2  = 44469  # D60  Fecha
$ws.Cells.Item(60, 10).Value2 = 500    # J60  Volumen
$ws.Cells.Item(60, 11).Value2 = 1200   # K60  Precio minimo
$ws.Cells.Item(60, 12).Value2 = 1200   # L60  Precio maximo
$ws.Cells.Item(60, 13).Value2 = 1200   # M60  Precio promedio ponderado
$ws.Cells.Item(60, 16).Value2 = 240    # P60  Precio $/Kg
